# Simulates a Power Query "Refresh" on the Query1 table that pulled in
# two new trial rows (COLO-PREVENT, TEST) and updated a couple of the
# "Days remaining" values, expanding the query table from A1:C10 to A1:C12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data that changed on refresh ---
$ws.Range("B8").Value = 15
$ws.Range("B10").Value = 35

# --- Append the two new rows pulled in by the refreshed query ---
$ws.Range("A11").Value = "COLO-PREVENT"
$ws.Range("C11").Value = 0
$ws.Range("A12").Value = "TEST"

# --- Power Query stamps the refreshed "Trial Name" column cells with an
#     explicit (General) number format, so match that for rows 2-12 ---
$ws.Range("A2:A12").NumberFormat = "General"

# --- Grow the query table / list object to cover the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C12"))

# --- Update the ExternalData_1 defined name range to match ---
$nm = $wb.Names.Item(1)
$nm.RefersTo = "=Sheet1!`$A`$1:`$C`$12"
